# Insert a new data row at worksheet row 579 (pushing the existing rows
# 579..662 down to 580..663), then populate the new row 579 with its
# values. Most fields duplicate what used to be in the old row 579
# (now shifted to row 580); only the Fecha/Volumen/Precio* columns
# (D, J, K, L, M, P) get new data values for the newly-inserted record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 579:662 down by one row, creating a blank row 579.
$ws.Rows.Item(579).Insert()

# Fill in the new row 579.
$ws.Cells.Item(579, 1).Value = 6
$ws.Cells.Item(579, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(579, 3).Value = "Metropolitana"
$ws.Cells.Item(579, 4).Value = 45154
$ws.Cells.Item(579, 5).Value = 13
$ws.Cells.Item(579, 6).Value = 100112043
$ws.Cells.Item(579, 7).Value = "Pepino ensalada"
$ws.Cells.Item(579, 8).Value = "Sin especificar"
$ws.Cells.Item(579, 9).Value = "Primera"
$ws.Cells.Item(579, 10).Value = 1100
$ws.Cells.Item(579, 11).Value = 8000
$ws.Cells.Item(579, 12).Value = 9000
$ws.Cells.Item(579, 13).Value = 8682
$ws.Cells.Item(579, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(579, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(579, 16).Value = 145
$ws.Cells.Item(579, 17).Value = 60
$ws.Cells.Item(579, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest
# of column D.
$ws.Cells.Item(579, 4).NumberFormat = $ws.Cells.Item(580, 4).NumberFormat
